# Fill in the previously-empty "Actual Result" (F) / "Test Result" (G) cells
# for TC_012 (favicon), TC_015 (hover effects), TC_016 (text contrast) and
# TC_007 (page load time), plus mark TC_014 / (the two NOT APPLICABLE rows)
# as not applicable.
#
# The order in which new text is first written matters: the workbook's
# shared-string table appends unique strings in first-use order, and the
# upstream edit appended them in the order favicon -> hover -> contrast ->
# page-load, so we replicate that order here even though the rows are
# touched out of top-to-bottom order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# TC_012 - "Verify the favicon is displayed" (row 13)
$ws.Range("F13").Value = "The correct favicon was displayed on the browser tab."
$ws.Range("F13").WrapText = $true
$ws.Range("G13").Value = "PASS"

# TC_015 - "Verify hover effects (if any)" (row 16)
$ws.Range("F16").Value = "Hover effects worked correctly, providing visual feedback when interacting with link elements."
$ws.Range("F16").WrapText = $true
$ws.Range("G16").Value = "PASS"
$ws.Rows(16).RowHeight = 72

# TC_016 - "Verify text contrast and readability" (row 17)
$ws.Range("F17").Value = "Text had sufficient contrast with the background for easy readability."
$ws.Range("F17").WrapText = $true
$ws.Range("G17").Value = "PASS"

# TC_007 - "Verify page performance and load time" (row 8)
$ws.Range("F8").Value = "The page load within a reasonable time (e.g., 2-3 seconds), without performance lag.                "
$ws.Range("F8").WrapText = $true
$ws.Range("G8").Value = "PASS"

# Rows 14 and 15 - mark as not applicable
$ws.Range("G14").Value = "NOT APPLICABLE"
$ws.Range("G15").Value = "NOT APPLICABLE"

# Restore the view to the top of the sheet and leave the cursor on G9,
# matching where the author's selection ended up after the edits.
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("G9").Select()
